# Fixed spelling errors in validation tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("A14").Value = '2. Warwick P A, D E Ramsey and M Takata. "Development of Hypertension from Unilateral Artery Stenosis in Conscious Dogs." Hypertension. 1990. Vol 16. No 4. p. 441-451.'

$ws = $wb.Worksheets.Item("HemorrhageClass2NoFluid")
$ws.Range("G1").Value = "Action Occurrence Time (s)"

$ws = $wb.Worksheets.Item("HemorrhageClass4NoFluid")
$ws.Range("G1").Value = "Action Occurrence Time (s)"

$ws = $wb.Worksheets.Item("HighAltitudeEnvironmentChange")
$ws.Range("G4").Value = "Action Occurrence Time (s)"

$ws = $wb.Worksheets.Item("WaterIngestion")
$ws.Range("G2").Value = "Action Occurrence Time (s)"

$ws = $wb.Worksheets.Item("Starvation")
$ws.Range("G2").Value = "Action Occurrence Time (s)"
$ws.Range("I2").Value = "Sampled Scenario Time (days)"
$ws.Range("C4").Value = "Patient is experiencing starvation"

$ws = $wb.Worksheets.Item("UnilateralStenosis")
$ws.Range("G1").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("BilateralStenosis")
$ws.Range("G1").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("HemorrhageClass2NoFluid")
$ws.Range("I1").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("HemorrhageClass4NoFluid")
$ws.Range("I1").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("HighAltitudeEnvironmentChange")
$ws.Range("I4").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("WaterIngestion")
$ws.Range("I2").Value = "Sampled Scenario Time (s)"

$ws = $wb.Worksheets.Item("BilateralStenosis")
$ws.Range("E3").Value = "90% bilateral occlusion of kidneys "

$ws = $wb.Worksheets.Item("UnilateralStenosis")
$ws.Range("E3").Value = "60% unilateral occlusion of kidneys "
